$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the existing D74:D81 plain formulas into an Excel "shared formula"
# group (t="shared", si="3"), matching how Excel behaves when the same
# formula is entered/filled across a contiguous range.
$ws.Range("D74:D81").Formula = "=C74/(24*60)"

# --- New row 82: 四方坪站 (site 四方坪站, date 2025-10-11 / serial 45941) ---
$ws.Range("A82").Value = 45941
$ws.Range("B82").Value = "四方坪站"
$ws.Range("C82").Formula = "=18946/126"
$ws.Range("D82").Formula = "=C82/(24*60)"
$ws.Range("E82").Formula = "=10508.12/126"
$ws.Range("F82").Formula = "=3652.19/126"
$ws.Range("G82").Formula = "=10508.12/(18946/60)"
$ws.Range("H82").Value = 3.5396825396825395

# --- New row 83: 高岭站 (site 高岭站, date 2025-10-11 / serial 45941) ---
$ws.Range("A83").Value = 45941
$ws.Range("B83").Value = "高岭站"
$ws.Range("C83").Formula = "=7009/36"
$ws.Range("D83").Formula = "=C83/(24*60)"
$ws.Range("E83").Formula = "=5191.68/36"
$ws.Range("F83").Formula = "=1371.67/36"
$ws.Range("G83").Formula = "=5191.68/(7009/60)"
$ws.Range("H83").Formula = "=189/36"

# Update the sheet's active selection to match the edited file (last cell
# touched, H83), mirroring the authored <selection activeCell="H83" .../>.
$ws.Range("H83").Select()
